# Bitacora_Laboratorio_Parcial_02_NRC_126.xlsx
# Ej. 2a. ley Newton - 9 de enero de 2024 - Lap HP
#
# - Gutierrez Garcia Diana Berenice (row 8): P3_Marco (I8) score 0 -> 5
#   (O8 = SUM(E8:N8) and P8 = (O8/50)*10 recalc automatically: 45->50, 9->10)
# - Turn on AutoFilter over the data range and register the hidden
#   _FilterDatabase defined name Excel creates for it.
# - Nudge the split/selection state of the sheet view towards the new
#   cursor position (bottom-left pane, cell B2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concentrado")
$ws.Activate() | Out-Null

# --- grade update -----------------------------------------------------
$ws.Range("I8").Value = 5

# --- AutoFilter ---------------------------------------------------------
$dataRange = $ws.Range("A1:Q16")
$dataRange.AutoFilter() | Out-Null

$filterName = $ws.Names.Add("_xlnm._FilterDatabase", $dataRange)
$filterName.Visible = $false

# --- sheet view / selection --------------------------------------------
# The sheet was already split between columns O/P (xSplit) and rows 1/2
# (ySplit); nudge the column split back to that same boundary (so the
# right-hand pane still starts at P1, like in the original file) and move
# the cursor to B2 (the new activeCell for the bottom-left pane).
$win = $excel.ActiveWindow
$win.SplitColumn = 15
$ws.Range("B2").Select() | Out-Null
